# Applies the refreshed crypto price/volume snapshot (GitHub Actions scrape).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.054.31"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").Value = "1.622.03"
$ws.Range("E3").Value = "  -0.90%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'213.88"
$ws.Range("E5").Value = "  -1.38%  "

$ws.Range("D6").Value = "'0.511"
$ws.Range("E6").Value = "  -1.29%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  +0.29%  "

$ws.Range("E9").Value = "  -1.53%  "

$ws.Range("D10").Value = "'19.90"
$ws.Range("E10").Value = "  -0.14%  "

$ws.Range("E11").Value = "  -0.62%  "

$ws.Range("D12").Value = "1.848.98"

$ws.Range("D13").Value = "1.625.53"
$ws.Range("E13").Value = "  -0.42%  "

$ws.Range("E14").Value = "  -0.25%  "

$ws.Range("D15").Value = "'0.538"
$ws.Range("E15").Value = "  -0.49%  "

$ws.Range("D16").Value = "27.048.21"
$ws.Range("E16").Value = "  -0.28%  "

$ws.Range("D17").Value = "'64.45"
$ws.Range("E17").Value = "  -3.29%  "

$ws.Range("D18").Value = "0.0₃0736"
$ws.Range("E18").Value = "  -0.42%  "

$ws.Range("D19").Value = "'214.21"
$ws.Range("E19").Value = "  -1.22%  "

$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").Value = "'6.82"
$ws.Range("E21").Value = "  -0.30%  "

$ws.Range("E22").Value = "  -1.84%  "

$ws.Range("E23").Value = "  -7.46%  "

$ws.Range("E24").Value = "  -1.06%  "

$ws.Range("D25").Value = "'147.33"
$ws.Range("E25").Value = "  +0.48%  "

$ws.Range("B26").Value = "BinanceUSD"
$ws.Range("C26").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'7.42"
$ws.Range("E27").Value = "  +0.33%  "

$ws.Range("E28").Value = "  -3.54%  "

$ws.Range("D29").Value = "'15.50"
$ws.Range("E29").Value = "  -1.01%  "

$ws.Range("E30").Value = "  +0.87%  "

$ws.Range("E31").Value = "  -1.13%  "

$ws.Range("E32").Value = "  -2.01%  "

$ws.Range("D33").Value = "'0.717"
$ws.Range("E33").Value = "  +31.97%  "

$ws.Range("D34").Value = "'3.00"
$ws.Range("E34").Value = "  -0.21%  "

$ws.Range("D35").Value = "1.334.50"
$ws.Range("E35").Value = "  +2.80%  "

$ws.Range("E36").Value = "  -1.09%  "

$ws.Range("E37").Value = "  -0.57%  "

$ws.Range("E38").Value = "  -0.34%  "

$ws.Range("D39").Value = "'0.838"
$ws.Range("E39").Value = "  -1.82%  "

$ws.Range("E40").Value = "  -0.05%  "

$ws.Range("E41").Value = "  -0.40%  "

$ws.Range("D42").Value = "'0.794"
$ws.Range("E42").Value = "  -1.51%  "

$ws.Range("D43").Value = "'5.35"
$ws.Range("E43").Value = "  +0.93%  "

$ws.Range("D44").Value = "'63.85"
$ws.Range("E44").Value = "  +3.57%  "

$ws.Range("D45").Value = "1.760.47"

$ws.Range("D46").Value = "'89.90"
$ws.Range("E46").Value = "  -1.57%  "

$ws.Range("E47").Value = "  +2.34%  "

$ws.Range("D48").Value = "'0.857"
$ws.Range("E48").Value = "  +28.31%  "

$ws.Range("E49").Value = "  -0.18%  "

$ws.Range("D50").Value = "'0.0994"
$ws.Range("E50").Value = "  +3.86%  "

$ws.Range("D51").Value = "'7.59"
$ws.Range("E51").Value = "  -0.77%  "

